$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B-E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values for columns B-E
$ws.Range("B2").Value = 0.82147052983003199
$ws.Range("C2").Value = 1.9278803310590156
$ws.Range("D2").Value = 1.1388024172627749
$ws.Range("E2").Value = 2.4494871755446774

# Update row 3 values for columns B-E
$ws.Range("B3").Value = 2.0547895786242916
$ws.Range("C3").Value = 5.9421813636307732
$ws.Range("D3").Value = 5.2553961548146289
$ws.Range("E3").Value = 3.154230747463485

# Update the selection to reflect the narrower range used in the edit
$ws.Range("B1:E3").Select() | Out-Null
